# Auto-generated: updates simulated-game transition matrix probabilities
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.1987179487179487
    "C2" = 0.5512820512820513
    "J2" = 0.01923076923076923
    "P2" = 0.1346153846153846
    "S2" = 0.09615384615384616
    "B3" = 0.003231017770597738
    "C3" = 0.02261712439418417
    "J3" = 0.04038772213247173
    "P3" = 0.7479806138933764
    "S3" = 0.1857835218093699
    "J4" = 0.0718562874251497
    "P4" = 0.6526946107784432
    "S4" = 0.2754491017964072
    "J5" = 0.25
    "P5" = 0.625
    "S5" = 0.125
    "B6" = 0.07407407407407407
    "D6" = 0.01371742112482853
    "F6" = 0.04938271604938271
    "J6" = 0.252400548696845
    "O6" = 0.01783264746227709
    "Q6" = 0.1467764060356653
    "R6" = 0.07270233196159122
    "S6" = 0.3731138545953361
    "B7" = 0.1143583227445997
    "D7" = 0.01905972045743329
    "E7" = 0.006353240152477764
    "F7" = 0.05844980940279543
    "J7" = 0.1423125794155019
    "O7" = 0.02287166454891995
    "Q7" = 0.1613722998729352
    "R7" = 0.06226175349428208
    "S7" = 0.4129606099110547
    "B8" = 0.09427414690572586
    "D8" = 0.02486986697513013
    "E8" = 0.002313475997686524
    "F8" = 0.05320994794679005
    "J8" = 0.1133603238866397
    "O8" = 0.01792943898207056
    "Q8" = 0.183342972816657
    "R8" = 0.09947946790052054
    "S8" = 0.4112203585887796
    "B9" = 0.09361069836552749
    "D9" = 0.02526002971768202
    "F9" = 0.0549777117384844
    "J9" = 0.1084695393759287
    "O9" = 0.01337295690936107
    "Q9" = 0.187221396731055
    "R9" = 0.09212481426448738
    "S9" = 0.424962852897474
    "B10" = 0.1102644554011654
    "D10" = 0.01927386822052891
    "E10" = 0.0006723442402510085
    "F10" = 0.06432093231734648
    "J10" = 0.1199013895114299
    "O10" = 0.01523980277902286
    "Q10" = 0.2140295831465711
    "R10" = 0.09188704616763783
    "S10" = 0.3644105782160466
    "F11" = 0.0008130081300813008
    "G11" = 0.1577235772357723
    "J11" = 0.09593495934959349
    "K11" = 0.2203252032520325
    "L11" = 0.5073170731707317
    "S11" = 0.01788617886178862
    "F12" = 0.001524390243902439
    "G12" = 0.7560975609756098
    "J12" = 0.1875
    "K12" = 0.01219512195121951
    "L12" = 0.01829268292682927
    "S12" = 0.02439024390243903
    "F13" = 0.005649717514124294
    "G13" = 0.6892655367231638
    "J13" = 0.231638418079096
    "S13" = 0.07344632768361582
    "F15" = 0.02309782608695652
    "H15" = 0.1861413043478261
    "I15" = 0.05706521739130434
    "J15" = 0.3451086956521739
    "K15" = 0.06793478260869565
    "M15" = 0.02445652173913044
    "O15" = 0.06793478260869565
    "S15" = 0.2282608695652174
    "F16" = 0.01147776183644189
    "H16" = 0.1865136298421808
    "I16" = 0.06743185078909612
    "J16" = 0.3830703012912482
    "K16" = 0.1248206599713056
    "M16" = 0.02439024390243903
    "O16" = 0.06025824964131994
    "S16" = 0.1420373027259684
    "F17" = 0.01724137931034483
    "H17" = 0.1724137931034483
    "I17" = 0.08312807881773399
    "J17" = 0.4113300492610837
    "K17" = 0.1028325123152709
    "M17" = 0.02032019704433497
    "N17" = 0.001231527093596059
    "O17" = 0.06342364532019705
    "S17" = 0.1280788177339902
    "F18" = 0.01196808510638298
    "H18" = 0.2446808510638298
    "I18" = 0.07446808510638298
    "J18" = 0.3962765957446808
    "K18" = 0.09707446808510638
    "M18" = 0.01462765957446809
    "O18" = 0.05319148936170213
    "S18" = 0.1077127659574468
    "F19" = 0.01769718155997378
    "H19" = 0.2173913043478261
    "I19" = 0.08542713567839195
    "J19" = 0.351321826523924
    "K19" = 0.1225693685820406
    "M19" = 0.02119292112737601
    "N19" = 0.0008739348918505571
    "O19" = 0.05964605636880053
    "S19" = 0.1238802709198165
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

Write-Host "Updated $($values.Count) cells"